$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hours for Week 3 (row 6) - Ryan's hours increased from 4 to 8
$ws.Range("B6").Value = 8

# Update hours for Week 4 (row 7) - new entries for Ryan and Tim
$ws.Range("B7").Value = 6
$ws.Range("D7").Value = 5

# Update hours for Week 5 (row 8) - new entries for Ryan and Tim
$ws.Range("B8").Value = 5
$ws.Range("D8").Value = 3

# Update hours for Week 6 (row 9) - new entry for Tim
$ws.Range("D9").Value = 4

# Update the active cell selection to B9
$ws.Range("B9").Select()
